$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "I", "Q", "R", "AC")

foreach ($col in $cols) {
    $rng7 = $ws.Range($col + "7")
    $rng8 = $ws.Range($col + "8")
    $v7 = $rng7.Value2
    $v8 = $rng8.Value2
    $rng7.Value = $v8
    $rng8.Value = $v7
}
